# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy the formatting used by the other header cells (bold,
# centered, bordered) from G1, then set the label.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data cells H2:H58: the per-row "Save" flag (0/1).
$saveVals = @(
    0,1,1,1,0,0,0,0,0,0,0,1,0,1,1,1,1,1,1,0,1,0,0,1,1,1,1,1,1,0,0,0,0,
    0,0,0,0,0,1,0,1,0,1,0,1,0,1,1,0,0,1,0,1,0,1,0,0
)

$arr = New-Object 'object[,]' $saveVals.Count,1
for ($i = 0; $i -lt $saveVals.Count; $i++) {
    $arr[$i,0] = $saveVals[$i]
}
$ws.Range("H2:H58").Value = $arr
